# The commit swaps the presentation's design from the custom "Integral"
# theme back to the built-in default "Office Theme" (colour scheme).
#
# Internally PowerPoint keeps one live/rendering theme (the one linked
# from the slide master - ppt/theme/theme2.xml in this deck) that every
# slide inherits its 12 theme colours from (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink). Re-applying the default Office colour
# scheme just rewrites those 12 colours; we do that through the
# ThemeColorScheme collection, which is exposed per-slide (and mirrors
# across every slide/layout because they all use <a:masterClrMapping/>).
#
# RGB() isn't available in this host, and the .RGB setter expects a
# 0xBBGGRR integer (standard COM colour order), so each target RRGGBB
# hex value below is written with its bytes reversed.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Colors(1).RGB  = 0x000000   # dk1      - 000000
$tcs.Colors(2).RGB  = 0xFFFFFF   # lt1      - FFFFFF
$tcs.Colors(3).RGB  = 0x6A5444   # dk2      - 44546A
$tcs.Colors(4).RGB  = 0xE6E6E7   # lt2      - E7E6E6
$tcs.Colors(5).RGB  = 0xD59B5B   # accent1  - 5B9BD5
$tcs.Colors(6).RGB  = 0x317DED   # accent2  - ED7D31
$tcs.Colors(7).RGB  = 0xA5A5A5   # accent3  - A5A5A5
$tcs.Colors(8).RGB  = 0x00C0FF   # accent4  - FFC000
$tcs.Colors(9).RGB  = 0xC47244   # accent5  - 4472C4
$tcs.Colors(10).RGB = 0x47AD70   # accent6  - 70AD47
$tcs.Colors(11).RGB = 0xC16305   # hlink    - 0563C1
$tcs.Colors(12).RGB = 0x724F95   # folHlink - 954F72
